# Weekly update: insert a new week of "Zanahoria" price records at the
# top of the data block (row 532) for Vega Central Mapocho de Santiago,
# pushing all the existing historical rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 532 (existing rows 532:643 shift to 536:647)
$ws.Rows("532:535").Insert()

# New week date: 2021-10-07 (serial 44476)
$newDate = Get-Date -Year 2021 -Month 10 -Day 7 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Row 532: Primera / Región Metropolitana
$ws.Cells.Item(532, 1).Value = 9
$ws.Cells.Item(532, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(532, 3).Value = "Metropolitana"
$ws.Cells.Item(532, 4).Value = $newDate
$ws.Cells.Item(532, 5).Value = 13
$ws.Cells.Item(532, 6).Value = 100114013
$ws.Cells.Item(532, 7).Value = "Zanahoria"
$ws.Cells.Item(532, 8).Value = "Sin especificar"
$ws.Cells.Item(532, 9).Value = "Primera"
$ws.Cells.Item(532, 10).Value = 160
$ws.Cells.Item(532, 11).Value = 8000
$ws.Cells.Item(532, 12).Value = 9000
$ws.Cells.Item(532, 13).Value = 8500
$ws.Cells.Item(532, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(532, 15).Value = "Región Metropolitana"
$ws.Cells.Item(532, 16).Value = 425
$ws.Cells.Item(532, 17).Value = 20
$ws.Cells.Item(532, 18).Value = "Hortaliza"

# Row 533: Primera / Región de La Araucanía
$ws.Cells.Item(533, 1).Value = 9
$ws.Cells.Item(533, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(533, 3).Value = "Metropolitana"
$ws.Cells.Item(533, 4).Value = $newDate
$ws.Cells.Item(533, 5).Value = 13
$ws.Cells.Item(533, 6).Value = 100114013
$ws.Cells.Item(533, 7).Value = "Zanahoria"
$ws.Cells.Item(533, 8).Value = "Sin especificar"
$ws.Cells.Item(533, 9).Value = "Primera"
$ws.Cells.Item(533, 10).Value = 187
$ws.Cells.Item(533, 11).Value = 8000
$ws.Cells.Item(533, 12).Value = 9000
$ws.Cells.Item(533, 13).Value = 8503
$ws.Cells.Item(533, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(533, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(533, 16).Value = 425
$ws.Cells.Item(533, 17).Value = 20
$ws.Cells.Item(533, 18).Value = "Hortaliza"

# Row 534: Segunda / Región Metropolitana
$ws.Cells.Item(534, 1).Value = 9
$ws.Cells.Item(534, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(534, 3).Value = "Metropolitana"
$ws.Cells.Item(534, 4).Value = $newDate
$ws.Cells.Item(534, 5).Value = 13
$ws.Cells.Item(534, 6).Value = 100114013
$ws.Cells.Item(534, 7).Value = "Zanahoria"
$ws.Cells.Item(534, 8).Value = "Sin especificar"
$ws.Cells.Item(534, 9).Value = "Segunda"
$ws.Cells.Item(534, 10).Value = 97
$ws.Cells.Item(534, 11).Value = 6000
$ws.Cells.Item(534, 12).Value = 7000
$ws.Cells.Item(534, 13).Value = 6495
$ws.Cells.Item(534, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(534, 15).Value = "Región Metropolitana"
$ws.Cells.Item(534, 16).Value = 325
$ws.Cells.Item(534, 17).Value = 20
$ws.Cells.Item(534, 18).Value = "Hortaliza"

# Row 535: Segunda / Región de La Araucanía
$ws.Cells.Item(535, 1).Value = 9
$ws.Cells.Item(535, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(535, 3).Value = "Metropolitana"
$ws.Cells.Item(535, 4).Value = $newDate
$ws.Cells.Item(535, 5).Value = 13
$ws.Cells.Item(535, 6).Value = 100114013
$ws.Cells.Item(535, 7).Value = "Zanahoria"
$ws.Cells.Item(535, 8).Value = "Sin especificar"
$ws.Cells.Item(535, 9).Value = "Segunda"
$ws.Cells.Item(535, 10).Value = 88
$ws.Cells.Item(535, 11).Value = 6000
$ws.Cells.Item(535, 12).Value = 7000
$ws.Cells.Item(535, 13).Value = 6500
$ws.Cells.Item(535, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(535, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(535, 16).Value = 325
$ws.Cells.Item(535, 17).Value = 20
$ws.Cells.Item(535, 18).Value = "Hortaliza"
